# Update "想去人数" (F column) counts across sheets, per commit
# "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 31
$ws1.Range("F4").Value = 191
$ws1.Range("F5").Value = 1097
$ws1.Range("F6").Value = 8244
$ws1.Range("F7").Value = 8244
$ws1.Range("F10").Value = 6913
$ws1.Range("F11").Value = 175
$ws1.Range("F12").Value = 5037
$ws1.Range("F13").Value = 5496
$ws1.Range("F14").Value = 1076
$ws1.Range("F15").Value = 335
$ws1.Range("F16").Value = 346
$ws1.Range("F18").Value = 315
$ws1.Range("F24").Value = 98
$ws1.Range("F25").Value = 9252
$ws1.Range("F26").Value = 73
$ws1.Range("F27").Value = 1685
$ws1.Range("F28").Value = 903
$ws1.Range("F31").Value = 1880
$ws1.Range("F37").Value = 1895
$ws1.Range("F38").Value = 243
$ws1.Range("F39").Value = 1208
$ws1.Range("F41").Value = 4834
$ws1.Range("F46").Value = 151
$ws1.Range("F47").Value = 1080
$ws1.Range("F49").Value = 923
$ws1.Range("F50").Value = 1267

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 34
$ws2.Range("F9").Value = 182

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 31
$ws4.Range("F4").Value = 191
$ws4.Range("F6").Value = 1097
$ws4.Range("F7").Value = 8244
$ws4.Range("F10").Value = 6913
$ws4.Range("F11").Value = 175
$ws4.Range("F14").Value = 5037
$ws4.Range("F15").Value = 5496
$ws4.Range("F16").Value = 1076
$ws4.Range("F17").Value = 335
$ws4.Range("F18").Value = 346
$ws4.Range("F20").Value = 315
$ws4.Range("F24").Value = 98
$ws4.Range("F25").Value = 9252
$ws4.Range("F26").Value = 73
$ws4.Range("F27").Value = 1685
$ws4.Range("F28").Value = 903
$ws4.Range("F31").Value = 1880
$ws4.Range("F37").Value = 1895
$ws4.Range("F38").Value = 243
$ws4.Range("F39").Value = 1208
$ws4.Range("F41").Value = 4834
$ws4.Range("F46").Value = 151
$ws4.Range("F47").Value = 1080
$ws4.Range("F49").Value = 923
$ws4.Range("F50").Value = 1267

